# Mise a jour de certains champs de Modules et de Professeurs
#
# Adds a new "Matières enseignés" column (E) to the professeur sheet,
# re-sizes the newly touched columns (C, D, E) and updates the active
# cell selection, mirroring the authored workbook change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell (also grows the shared-string table and the sheet
# dimension to A1:E1 automatically).
$ws.Range("E1").Value = "Matières enseignés"

# Custom column widths for C, D and E (values as shown in the Excel
# "Column Width" dialog, i.e. the ColumnWidth COM property).
$ws.Columns.Item(3).ColumnWidth = 26.736979166666668
$ws.Columns.Item(4).ColumnWidth = 14.877604166666666
$ws.Columns.Item(5).ColumnWidth = 30.877604166666668

# Move the active selection to E6, as left by the author.
[void]$ws.Range("E6").Select()
